$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.372.38"
$ws.Range("E2").Value = "  -2.94%  "

$ws.Range("D3").Value = "1.952.19"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("E4").Value = "  -1.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.66%  "

$ws.Range("E6").Value = "  -1.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4765"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4021"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08429"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.052"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.10%  "

$ws.Range("D13").Value = "1.955.32"
$ws.Range("E13").Value = "  -2.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.537"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.138"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001068"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06582"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.49%  "

$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.812"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.83%  "

$ws.Range("D23").Value = "28.389.18"

$ws.Range("E24").Value = "  -3.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.280"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("D26").Value = "2.186.77"
$ws.Range("E26").Value = "  -4.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.908"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.58%  "

$ws.Range("E30").Value = "  -6.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9760"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09608"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.448"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.589"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.658"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.956"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02323"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06206"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.248"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6195"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.21%  "

$ws.Range("E42").Value = "  -3.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.005"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1913"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.344"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5946"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.59%  "

$ws.Range("E48").Value = "  -5.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.383"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000317"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06801"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.95%  "
